# Updates cryptos list values (Price / Volume(1h), and for two swapped
# ranking pairs also Coin name + Link) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.247.92'
$ws.Range("E2").Value = '  -0.93%  '

# Row 3
$ws.Range("D3").Value = '2.247.00'
$ws.Range("E3").Value = '  -0.91%  '

# Row 4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").Value = '''247.06'

# Row 6
$ws.Range("D6").Value = '''0.621'
$ws.Range("E6").Value = '  -3.32%  '

# Row 7
$ws.Range("D7").Value = '''74.57'
$ws.Range("E7").Value = '  -1.35%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").Value = '''0.615'
$ws.Range("E9").Value = '  -4.75%  '

# Row 10
$ws.Range("D10").Value = '''42.50'
$ws.Range("E10").Value = '  +6.31%  '

# Row 11
$ws.Range("D11").Value = '''0.0946'
$ws.Range("E11").Value = '  -2.48%  '

# Row 12
$ws.Range("E12").Value = '  -2.25%  '

# Row 13
$ws.Range("E13").Value = '  -3.38%  '

# Row 14
$ws.Range("E14").Value = '  -2.88%  '

# Row 15
$ws.Range("D15").Value = '''0.853'
$ws.Range("E15").Value = '  -1.46%  '

# Row 16
$ws.Range("D16").Value = '2.246.33'
$ws.Range("E16").Value = '  -0.89%  '

# Row 17
$ws.Range("D17").Value = '42.110.64'
$ws.Range("E17").Value = '  -1.02%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0987'
$ws.Range("E18").Value = '  -0.61%  '

# Row 19
$ws.Range("E19").Value = '  -0.25%  '

# Row 20
$ws.Range("D20").Value = '''71.96'
$ws.Range("E20").Value = '  -0.09%  '

# Row 21
$ws.Range("E21").Value = '  +3.80%  '

# Row 22
$ws.Range("D22").Value = '''230.44'
$ws.Range("E22").Value = '  -2.24%  '

# Row 23
$ws.Range("D23").Value = '''8.94'
$ws.Range("E23").Value = '  +38.23%  '

# Row 24
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("D25").Value = '''11.28'
$ws.Range("E25").Value = '  +0.46%  '

# Row 26
$ws.Range("D26").Value = '''3.63'
$ws.Range("E26").Value = '  -5.18%  '

# Row 27
$ws.Range("D27").Value = '''2.31'
$ws.Range("E27").Value = '  -2.81%  '

# Row 28
$ws.Range("D28").Value = '''2.23'
$ws.Range("E28").Value = '  +4.38%  '

# Row 29
$ws.Range("D29").Value = '''169.24'
$ws.Range("E29").Value = '  +0.87%  '

# Row 30
$ws.Range("D30").Value = '''20.71'
$ws.Range("E30").Value = '  -0.88%  '

# Row 31
$ws.Range("D31").Value = '''0.0825'
$ws.Range("E31").Value = '  -3.28%  '

# Row 32
$ws.Range("E32").Value = '  -3.82%  '

# Row 33
$ws.Range("D33").Value = '''30.34'
$ws.Range("E33").Value = '  -5.68%  '

# Row 34
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '''0.125'
$ws.Range("E34").Value = '  -1.45%  '

# Row 35
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '''5.20'
$ws.Range("E35").Value = '  +10.04%  '

# Row 36
$ws.Range("E36").Value = '  +0.16%  '

# Row 37
$ws.Range("E37").Value = '  -0.39%  '

# Row 38
$ws.Range("E38").Value = '  -0.63%  '

# Row 39
$ws.Range("D39").Value = '''2.19'
$ws.Range("E39").Value = '  -2.97%  '

# Row 40
$ws.Range("E40").Value = '  -0.76%  '

# Row 41
$ws.Range("D41").Value = '''61.87'
$ws.Range("E41").Value = '  +1.10%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.203'
$ws.Range("E42").Value = '  -1.72%  '

# Row 43
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''108.73'
$ws.Range("E43").Value = '  +2.37%  '

# Row 44
$ws.Range("E44").Value = '  -2.09%  '

# Row 45
$ws.Range("E45").Value = '  +1.40%  '

# Row 46
$ws.Range("D46").Value = '''0.996'
$ws.Range("E46").Value = '  -0.39%  '

# Row 47
$ws.Range("E47").Value = '  -2.74%  '

# Row 48
$ws.Range("E48").Value = '  -0.37%  '

# Row 49
$ws.Range("D49").Value = '''2.30'
$ws.Range("E49").Value = '  +2.74%  '

# Row 50
$ws.Range("D50").Value = '''4.18'
$ws.Range("E50").Value = '  -11.38%  '

# Row 51
$ws.Range("D51").Value = '''4.12'
$ws.Range("E51").Value = '  -1.58%  '
